$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update test case name strings (B6, B7, B9, B12)
$ws.Range("B6").Value = "TC_PM_APFA_ListView_D2"
$ws.Range("B7").Value = "TC_PM_APFA_DebtInvestmentProfile_AddNewProfile"
$ws.Range("B9").Value = "TC_PM_APFA_DebtInvestmentProfile_EditProfile"
$ws.Range("B12").Value = "TC_PM_APFA_DebtInvestmentProfile_DeleteProfile"

# Update sheet view: move selection to B6 (also resets the scrolled topLeftCell back to default)
$ws.Activate()
$ws.Range("B6").Select()
